$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking
# strings like "1.015" or "28.432.64" are not coerced into numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.432.64"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "1.878.98"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("D4").Value = "1.015"
$ws.Range("E4").Value = "  +0.96%  "
$ws.Range("D5").Value = "315.91"
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("D6").Value = "1.013"
$ws.Range("E6").Value = "  +1.00%  "
$ws.Range("D7").Value = "0.5139"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "0.3947"
$ws.Range("E8").Value = "  +1.27%  "
$ws.Range("D9").Value = "0.08337"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("D10").Value = "1.122"
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("D11").Value = "41.96"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").Value = "6.277"
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("D13").Value = "1.876.79"
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("D14").Value = "20.42"
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("D15").Value = "7.264"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "1.014"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").Value = "0.00001107"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("D19").Value = "0.06730"
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("D20").Value = "17.80"
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").Value = "1.013"
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").Value = "5.988"
$ws.Range("E22").Value = "  -0.96%  "
$ws.Range("D23").Value = "28.469.20"
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("D24").Value = "11.17"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "2.250"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("D26").Value = "2.085.94"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").Value = "161.55"
$ws.Range("E27").Value = "  +1.98%  "
$ws.Range("D28").Value = "20.69"
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("D29").Value = "2.426"
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("D30").Value = "127.42"
$ws.Range("E30").Value = "  +1.78%  "
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").Value = "1.045"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").Value = "5.882"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("D34").Value = "3.634"
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("E35").Value = "  +0.48%  "
$ws.Range("D36").Value = "0.06529"
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").Value = "9.203"
$ws.Range("E37").Value = "  -5.42%  "
$ws.Range("D38").Value = "0.2190"
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("D39").Value = "1.255"
$ws.Range("E39").Value = "  +2.12%  "
$ws.Range("D40").Value = "0.6479"
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("D41").Value = "1.191"
$ws.Range("E41").Value = "  -1.33%  "
$ws.Range("D42").Value = "5.006"
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("D43").Value = "11.15"
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("D44").Value = "0.6066"
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("D45").Value = "13.13"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("D46").Value = "3.700"
$ws.Range("E46").Value = "  +0.66%  "
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("D48").Value = "2.015"
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("D49").Value = "1.214"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("D50").Value = "121.76"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").Value = "0.06900"
$ws.Range("E51").Value = "  +0.23%  "

# Restore the original (default) cell style on column D now that the
# text values are safely stored, so formatting matches the source file.
$priceRange.Style = "Normal"
